$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1), columns K..T -----------------------------
# Copy the header formatting (bold/centered/bordered style used by A1:J1)
# onto the new header cells, then fill in the labels.
$ws.Range("A1").Copy()
$ws.Range("K1:T1").PasteSpecial(-4122)

$ws.Range("K1").Value = "Nome da Ferramenta 1"
$ws.Range("L1").Value = "Objetivo 1"
$ws.Range("M1").Value = "Categoria da Ferramenta 1"
$ws.Range("N1").Value = "Importância 1"
$ws.Range("O1").Value = "Horas Gastas Mensais 1"
$ws.Range("P1").Value = "Nome da Ferramenta 2"
$ws.Range("Q1").Value = "Objetivo 2"
$ws.Range("R1").Value = "Categoria da Ferramenta 2"
$ws.Range("S1").Value = "Importância 2"
$ws.Range("T1").Value = "Horas Gastas Mensais 2"

# --- Back-fill existing rows 2-10 with blank cells in K..T ---------------
# Copy an already-blank data cell's (unstyled) formatting across the new
# columns so the cells exist (matching the widened used-range) without
# picking up the header style.
$ws.Range("F2").Copy()
$ws.Range("K2:T10").PasteSpecial(-4122)

# --- New row 11 ------------------------------------------------------------
$ws.Range("A11").Value = "pedro.paulistano@mrv.com.br"
$ws.Range("F11").Value = "Painel do Portifólio - Planejamento da Produção - PLNESROBR004"
$ws.Range("I11").Value = "2025-05-19 20:06:51"
$ws.Range("J11").Value = "Painel do Portifólio - Planejamento da Produção - PLNESROBR004: asdfghjklç`n"

$ws.Range("K11").Value = "Planilha automatizada"
$ws.Range("L11").Value = "ok"
$ws.Range("M11").Value = "Painel Power BI"
$ws.Range("N11").Value = 1
$ws.Range("O11").NumberFormat = "@"
$ws.Range("O11").Value = "15"
$ws.Range("P11").Value = "Planilha geral - teste"
$ws.Range("Q11").Value = "Teste Geração - Talita 03.04.25 v03"
$ws.Range("R11").Value = "Painel Power BI"
$ws.Range("S11").Value = 4
$ws.Range("T11").NumberFormat = "@"
$ws.Range("T11").Value = "21"

# J11 contains an embedded newline, which otherwise causes the row to pick
# up an explicit (and unwanted) custom height; AutoFit restores the default.
$ws.Rows("11").AutoFit()
